$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "1811"
$ws.Range("E17").Value = "1812"
$ws.Range("E18").Value = "1901"
$ws.Range("E19").Value = "1902"
$ws.Range("E20").Value = "1903"

$ws.Range("F16").Value = 31249
$ws.Range("F20").Value = 26041
